$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 46: was "lowest common ancestor of a binary tree" -> now "time based key-value store"
$ws.Range("A46").Value = "time based key-value store"
$ws.Range("F46").Value = "https://leetcode.com/problems/time-based-key-value-store/"
$ws.Range("B46").Value = "bisect"
$ws.Range("C46").Value = "dict of dicts"
$ws.Range("E46").Value = 36

# --- Row 47 (new): accounts merge
$ws.Range("F47").Value = "https://leetcode.com/problems/accounts-merge/"
$ws.Range("A47").Value = "accounts merge"
$ws.Range("B47").Value = "graph"
$ws.Range("C47").Value = "adjacency table"
$ws.Range("D47").Value = "dfs"
$ws.Range("E47").Value = 90

# --- Row 48 (new): sort colors
$ws.Range("F48").Value = "https://leetcode.com/problems/sort-colors/"
$ws.Range("A48").Value = "sort colors"
$ws.Range("B48").Value = "counter"
$ws.Range("C48").Value = "array"
$ws.Range("E48").Value = 13

# F47/F48 use the same "URL text" direct formatting (Consolas green, vertical center)
# already applied to F1 (style index 1 in the original workbook) - copy that formatting over.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("F47").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Copy() | Out-Null
$ws.Range("F48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Match final selection state
$ws.Range("D48").Select() | Out-Null
